# Apply metadata.xlsx edits: add new columns of descriptive data to row 1 (headers)
# and row 2 (values), widen columns F/G, normalize the font used by style index 1
# back to the workbook default, and update the active view/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header (row 1) ---
$ws.Range("U1").Value = "Subject"

# --- New data (row 2) ---
$ws.Range("G2").Value = "Arnold and Deanne Kaplan Collection of Early American Judaica (University of Pennsylvania)"
$ws.Range("H2").Value = "Arc.MS.56"
$ws.Range("I2").Value = "Trade cards"
$ws.Range("K2").Value = "English"
$ws.Range("L2").Value = "undated"
$ws.Range("N2").Value = "J. Rosenblatt & Co."
$ws.Range("O2").Value = "Baltimore, Maryland, United States | Maryland, United States"
$ws.Range("P2").Value = "J. Rosenblatt & Co.: Importers: Earthenware, China, Majolica, Novelties | 32 South Howard Street, Baltimore, MD"
$ws.Range("T2").Value = "http://rightsstatements.org/page/NoC-US/1.0/?"
$ws.Range("U2").Value = "House furnishings | Jewish merchants | Trade cards (advertising)"

# --- Column widths ---
$ws.Columns.Item(6).ColumnWidth = 82.83203125
$ws.Columns.Item(7).ColumnWidth = 18

# --- Reset the font used across row 1 & 2 back to the workbook default font ---
$headerRow = $ws.Range("A1:U2")
$headerRow.Font.Name = "Calibri"
$headerRow.Font.Size = 12
$headerRow.Font.Color = $ws.Range("A2").Font.Color
$headerRow.Font.ThemeColor = 1

# --- Selection / view ---
$ws.Range("G2").Select()
$ws.Rows.Item(2).EntireRow.Select()
$excel.ActiveWindow.ScrollColumn = 7
